# "Generate Report for Handoff"
#
# The localization-status report moves from "Handed back: in sync with
# en-US" to "Ready for handoff", and the associated timestamps advance a
# few seconds/minutes to the moment the handoff report was (re)generated.
# A couple of the "Status" columns also get narrower (their old width was
# sized for the long "Handed back..." text).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Timestamps bumped forward to the new handoff-report generation time ---
$wsOverview.Range("G2").Value = "2016-08-16 22:58:00"
$wsDeDe.Range("H2").Value = "2016-08-16 22:58:00"
$wsZhCn.Range("H2").Value = "2016-08-16 22:57:54"

# --- Narrow the (now shorter) Status columns ---
# ColumnWidth is quantized by the host to the nearest 1/6 character, so we
# feed it a value that snaps to the closest reachable width to the target.
$newColWidth = 16.333333333333336

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
